$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 899.6667
$ws.Cells.Item(19, 9).Value = 899.6667
$ws.Cells.Item(19, 11).Value = 899.6667
$ws.Cells.Item(19, 13).Value = -724.6667
$ws.Cells.Item(33, 8).Value = 8640758
$ws.Cells.Item(33, 9).Value = 11724951
$ws.Cells.Item(33, 11).Value = 11724951
$ws.Cells.Item(33, 13).Value = -11724722
$ws.Cells.Item(62, 8).Value = 4450.8335
$ws.Cells.Item(62, 10).Value = 4004
$ws.Cells.Item(62, 12).Value = 4004
$ws.Cells.Item(62, 14).Value = -5252
$ws.Cells.Item(65, 8).Value = 4450.8335
$ws.Cells.Item(65, 10).Value = 4004
$ws.Cells.Item(65, 12).Value = 20020
$ws.Cells.Item(65, 14).Value = -26260
$ws.Cells.Item(98, 8).Value = 3988.9333
$ws.Cells.Item(98, 9).Value = 3023.625
$ws.Cells.Item(98, 11).Value = 3023.625
$ws.Cells.Item(98, 13).Value = -1525.625
$ws.Cells.Item(122, 8).Value = 3988.9333
$ws.Cells.Item(122, 9).Value = 3023.625
$ws.Cells.Item(122, 11).Value = 9070.875
$ws.Cells.Item(122, 13).Value = -6620.875
$ws.Cells.Item(125, 8).Value = 202272860
$ws.Cells.Item(125, 9).Value = 337121200
$ws.Cells.Item(125, 10).Value = 335
$ws.Cells.Item(125, 11).Value = 3034090800
$ws.Cells.Item(125, 12).Value = 3015
$ws.Cells.Item(125, 13).Value = -3034088340
$ws.Cells.Item(125, 14).Value = -7935
$ws.Cells.Item(132, 8).Value = 1334.0938
$ws.Cells.Item(132, 9).Value = 1162.138
$ws.Cells.Item(132, 10).Value = 2996.3333
$ws.Cells.Item(132, 11).Value = 3486.414
$ws.Cells.Item(132, 12).Value = 8988.999899999999
$ws.Cells.Item(132, 13).Value = -956.4139999999998
$ws.Cells.Item(132, 14).Value = -14048.9999
$ws.Cells.Item(135, 8).Value = 1426.5
$ws.Cells.Item(135, 9).Value = 1108.2
$ws.Cells.Item(135, 11).Value = 9973.800000000001
$ws.Cells.Item(135, 13).Value = -7438.800000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 489.22726
$ws.Cells.Item(5, 9).Value = 545.3570999999999
$ws.Cells.Item(5, 10).Value = 391
$ws.Cells.Item(5, 11).Value = 545.3570999999999
$ws.Cells.Item(5, 12).Value = 391
$ws.Cells.Item(5, 13).Value = -433.3570999999999
$ws.Cells.Item(5, 14).Value = -615
$ws.Cells.Item(32, 8).Value = 210895.48
$ws.Cells.Item(32, 9).Value = 244939.56
$ws.Cells.Item(32, 10).Value = 11494.429
$ws.Cells.Item(32, 11).Value = 244939.56
$ws.Cells.Item(32, 12).Value = 11494.429
$ws.Cells.Item(32, 13).Value = -244652.56
$ws.Cells.Item(32, 14).Value = -12068.429
$ws.Cells.Item(45, 8).Value = 50566.43
$ws.Cells.Item(45, 9).Value = 69154.47
$ws.Cells.Item(45, 11).Value = 69154.47
$ws.Cells.Item(45, 13).Value = -68777.47
$ws.Cells.Item(132, 8).Value = 1922.2821
$ws.Cells.Item(132, 9).Value = 1760.9584
$ws.Cells.Item(132, 10).Value = 3858.1667
$ws.Cells.Item(132, 11).Value = 5282.8752
$ws.Cells.Item(132, 12).Value = 11574.5001
$ws.Cells.Item(132, 13).Value = -2752.8752
$ws.Cells.Item(132, 14).Value = -16634.5001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 489.22726
$ws.Cells.Item(4, 9).Value = 545.3570999999999
$ws.Cells.Item(4, 10).Value = 391
$ws.Cells.Item(4, 11).Value = 545.3570999999999
$ws.Cells.Item(4, 12).Value = 391
$ws.Cells.Item(4, 13).Value = -430.3570999999999
$ws.Cells.Item(4, 14).Value = -621
$ws.Cells.Item(80, 8).Value = 100000260
$ws.Cells.Item(80, 9).Value = 250000050
$ws.Cells.Item(80, 11).Value = 250000050
$ws.Cells.Item(80, 13).Value = -249999052
$ws.Cells.Item(83, 8).Value = 100000260
$ws.Cells.Item(83, 9).Value = 250000050
$ws.Cells.Item(83, 11).Value = 1250000250
$ws.Cells.Item(83, 13).Value = -1249995258
$ws.Cells.Item(107, 8).Value = 7925.4165
$ws.Cells.Item(107, 9).Value = 9199.138000000001
$ws.Cells.Item(107, 10).Value = 2648.5715
$ws.Cells.Item(107, 11).Value = 9199.138000000001
$ws.Cells.Item(107, 12).Value = 2648.5715
$ws.Cells.Item(107, 13).Value = -7279.138000000001
$ws.Cells.Item(107, 14).Value = -6488.5715
$ws.Cells.Item(134, 8).Value = 18001650
$ws.Cells.Item(134, 9).Value = 1317.9762
$ws.Cells.Item(134, 10).Value = 112503400
$ws.Cells.Item(134, 11).Value = 3953.9286
$ws.Cells.Item(134, 12).Value = 337510200
$ws.Cells.Item(134, 13).Value = -1418.9286
$ws.Cells.Item(134, 14).Value = -337515270

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 97.7619
$ws.Cells.Item(7, 9).Value = 101.4
$ws.Cells.Item(7, 11).Value = 101.4
$ws.Cells.Item(7, 13).Value = 11.59999999999999
$ws.Cells.Item(31, 8).Value = 2588.2263
$ws.Cells.Item(31, 9).Value = 2172.0571
$ws.Cells.Item(31, 11).Value = 2172.0571
$ws.Cells.Item(31, 13).Value = -1877.0571
$ws.Cells.Item(34, 8).Value = 2588.2263
$ws.Cells.Item(34, 9).Value = 2172.0571
$ws.Cells.Item(34, 11).Value = 2172.0571
$ws.Cells.Item(34, 13).Value = -1970.0571
$ws.Cells.Item(58, 8).Value = 1519.4728
$ws.Cells.Item(58, 9).Value = 1473.4865
$ws.Cells.Item(58, 10).Value = 1614
$ws.Cells.Item(58, 11).Value = 1473.4865
$ws.Cells.Item(58, 12).Value = 1614
$ws.Cells.Item(58, 13).Value = -1270.4865
$ws.Cells.Item(58, 14).Value = -2020
$ws.Cells.Item(86, 8).Value = 18385.428
$ws.Cells.Item(86, 9).Value = 4759.6
$ws.Cells.Item(86, 10).Value = 52450
$ws.Cells.Item(86, 11).Value = 4759.6
$ws.Cells.Item(86, 12).Value = 52450
$ws.Cells.Item(86, 13).Value = -3636.6
$ws.Cells.Item(86, 14).Value = -54696
$ws.Cells.Item(89, 8).Value = 18385.428
$ws.Cells.Item(89, 9).Value = 4759.6
$ws.Cells.Item(89, 10).Value = 52450
$ws.Cells.Item(89, 11).Value = 23798
$ws.Cells.Item(89, 12).Value = 262250
$ws.Cells.Item(89, 13).Value = -18182
$ws.Cells.Item(89, 14).Value = -273482
$ws.Cells.Item(99, 8).Value = 104001360
$ws.Cells.Item(99, 9).Value = 10000750
$ws.Cells.Item(99, 11).Value = 10000750
$ws.Cells.Item(99, 13).Value = -9999252
$ws.Cells.Item(105, 8).Value = 2921.875
$ws.Cells.Item(105, 9).Value = 2008.4
$ws.Cells.Item(105, 10).Value = 4444.3335
$ws.Cells.Item(105, 11).Value = 2008.4
$ws.Cells.Item(105, 12).Value = 4444.3335
$ws.Cells.Item(105, 13).Value = -261.4000000000001
$ws.Cells.Item(105, 14).Value = -7938.3335
$ws.Cells.Item(107, 8).Value = 2188.7144
$ws.Cells.Item(107, 9).Value = 1797.3572
$ws.Cells.Item(107, 10).Value = 2971.4285
$ws.Cells.Item(107, 11).Value = 1797.3572
$ws.Cells.Item(107, 12).Value = 2971.4285
$ws.Cells.Item(107, 13).Value = 122.6428000000001
$ws.Cells.Item(107, 14).Value = -6811.4285
$ws.Cells.Item(126, 8).Value = 104001360
$ws.Cells.Item(126, 9).Value = 10000750
$ws.Cells.Item(126, 11).Value = 30002250
$ws.Cells.Item(126, 13).Value = -29999780
$ws.Cells.Item(132, 8).Value = 20571.057
$ws.Cells.Item(132, 9).Value = 23260.37
$ws.Cells.Item(132, 10).Value = 2898.4285
$ws.Cells.Item(132, 11).Value = 69781.11
$ws.Cells.Item(132, 12).Value = 8695.2855
$ws.Cells.Item(132, 13).Value = -67251.11
$ws.Cells.Item(132, 14).Value = -13755.2855
$ws.Cells.Item(134, 8).Value = 1469.6818
$ws.Cells.Item(134, 9).Value = 1290.8125
$ws.Cells.Item(134, 11).Value = 3872.4375
$ws.Cells.Item(134, 13).Value = -1337.4375
$ws.Cells.Item(136, 8).Value = 1519.4728
$ws.Cells.Item(136, 9).Value = 1473.4865
$ws.Cells.Item(136, 10).Value = 1614
$ws.Cells.Item(136, 11).Value = 4420.4595
$ws.Cells.Item(136, 12).Value = 4842
$ws.Cells.Item(136, 13).Value = -1870.4595
$ws.Cells.Item(136, 14).Value = -9942

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(101, 8).Value = 19496.25
$ws.Cells.Item(101, 10).Value = 19496.25
$ws.Cells.Item(101, 12).Value = 58488.75
$ws.Cells.Item(101, 14).Value = -63356.75
$ws.Cells.Item(114, 8).Value = 2015.6923
$ws.Cells.Item(114, 9).Value = 1643
$ws.Cells.Item(114, 10).Value = 2248.625
$ws.Cells.Item(114, 11).Value = 4929
$ws.Cells.Item(114, 12).Value = 6745.875
$ws.Cells.Item(114, 13).Value = -1675
$ws.Cells.Item(114, 14).Value = -13253.875
$ws.Cells.Item(131, 8).Value = 7579470
$ws.Cells.Item(131, 10).Value = 4000
$ws.Cells.Item(131, 12).Value = 12000
$ws.Cells.Item(131, 14).Value = -22080

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(21, 8).Value = 2015800
$ws.Cells.Item(21, 10).Value = 20000
$ws.Cells.Item(21, 12).Value = 20000
$ws.Cells.Item(21, 14).Value = -20346
$ws.Cells.Item(24, 8).Value = 2509224.8
$ws.Cells.Item(24, 9).Value = 5006700
$ws.Cells.Item(24, 11).Value = 5006700
$ws.Cells.Item(24, 13).Value = -5006527
$ws.Cells.Item(30, 8).Value = 2015800
$ws.Cells.Item(30, 10).Value = 20000
$ws.Cells.Item(30, 12).Value = 20000
$ws.Cells.Item(30, 14).Value = -20210
$ws.Cells.Item(80, 8).Value = 21798474
$ws.Cells.Item(80, 9).Value = 95560.16
$ws.Cells.Item(80, 10).Value = 50012260
$ws.Cells.Item(80, 11).Value = 95560.16
$ws.Cells.Item(80, 12).Value = 50012260
$ws.Cells.Item(80, 13).Value = -94562.16
$ws.Cells.Item(80, 14).Value = -50014256
$ws.Cells.Item(83, 8).Value = 21798474
$ws.Cells.Item(83, 9).Value = 95560.16
$ws.Cells.Item(83, 10).Value = 50012260
$ws.Cells.Item(83, 11).Value = 477800.8
$ws.Cells.Item(83, 12).Value = 250061300
$ws.Cells.Item(83, 13).Value = -472808.8
$ws.Cells.Item(83, 14).Value = -250071284
$ws.Cells.Item(107, 8).Value = 77714.766
$ws.Cells.Item(107, 9).Value = 333601
$ws.Cells.Item(107, 11).Value = 333601
$ws.Cells.Item(107, 13).Value = -331681
$ws.Cells.Item(132, 8).Value = 737970.7
$ws.Cells.Item(132, 9).Value = 6154.625
$ws.Cells.Item(132, 10).Value = 2334660.2
$ws.Cells.Item(132, 11).Value = 18463.875
$ws.Cells.Item(132, 12).Value = 7003980.600000001
$ws.Cells.Item(132, 13).Value = -15933.875
$ws.Cells.Item(132, 14).Value = -7009040.600000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1065.3226
$ws.Cells.Item(16, 9).Value = 1052.8
$ws.Cells.Item(16, 10).Value = 1117.5
$ws.Cells.Item(16, 11).Value = 1052.8
$ws.Cells.Item(16, 12).Value = 1117.5
$ws.Cells.Item(16, 13).Value = -882.8
$ws.Cells.Item(16, 14).Value = -1457.5
$ws.Cells.Item(100, 8).Value = 3255.6667
$ws.Cells.Item(100, 9).Value = 3509.8
$ws.Cells.Item(100, 10).Value = 1985
$ws.Cells.Item(100, 11).Value = 3509.8
$ws.Cells.Item(100, 12).Value = 1985
$ws.Cells.Item(100, 13).Value = -2968.8
$ws.Cells.Item(100, 14).Value = -3067

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 893475.5600000001
$ws.Cells.Item(107, 9).Value = 636.8929000000001
$ws.Cells.Item(107, 10).Value = 7143346
$ws.Cells.Item(107, 11).Value = 1910.6787
$ws.Cells.Item(107, 12).Value = 21430038
$ws.Cells.Item(107, 13).Value = 9.321299999999837
$ws.Cells.Item(107, 14).Value = -21433878
$ws.Cells.Item(122, 8).Value = 1677.6
$ws.Cells.Item(122, 9).Value = 1425.5333
$ws.Cells.Item(122, 11).Value = 4276.5999
$ws.Cells.Item(122, 13).Value = -1826.5999
$ws.Cells.Item(132, 8).Value = 2227.6428
$ws.Cells.Item(132, 9).Value = 1699.3334
$ws.Cells.Item(132, 11).Value = 5098.0002
$ws.Cells.Item(132, 13).Value = -2568.0002
$ws.Cells.Item(136, 8).Value = 25445.5
$ws.Cells.Item(136, 9).Value = 32893.676
$ws.Cells.Item(136, 10).Value = 4455.1816
$ws.Cells.Item(136, 11).Value = 98681.02799999999
$ws.Cells.Item(136, 12).Value = 13365.5448
$ws.Cells.Item(136, 13).Value = -96131.02799999999
$ws.Cells.Item(136, 14).Value = -18465.5448
